$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swaps the full contents (A:AY) of two rows, using a scratch row as a buffer.
# NOTE: positional parameters only -- named parameters (e.g. "-RowA 12") are not
# reliably bound inside functions in this host, so call this with positional args.
function Swap-Rows {
    param($RowA, $RowB)

    $LastCol = "AY"
    $TempRow = 1000

    $rangeA = $ws.Range("A" + $RowA + ":" + $LastCol + $RowA)
    $rangeB = $ws.Range("A" + $RowB + ":" + $LastCol + $RowB)
    $rangeTemp = $ws.Range("A" + $TempRow + ":" + $LastCol + $TempRow)

    # Stash row A in the scratch row.
    $rangeTemp.ClearContents()
    $rangeA.Copy()
    $rangeTemp.PasteSpecial(-4104)  # xlPasteAll

    # Move row B's content into row A's slot (clear first so cells that are
    # blank in B don't keep stale values from A).
    $rangeA.ClearContents()
    $rangeB.Copy()
    $rangeA.PasteSpecial(-4104)

    # Move the stashed row A content into row B's slot.
    $rangeB.ClearContents()
    $rangeTemp.Copy()
    $rangeB.PasteSpecial(-4104)

    $rangeTemp.ClearContents()
    $excel.CutCopyMode = 0
}

Swap-Rows 12 14
Swap-Rows 16 17
Swap-Rows 23 24
